$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header in B1 from "pheno" to "phenotype"
$ws.Range("B1").Value = "phenotype"

# Mirror Excel's typical behavior of leaving the active cell on B2 after the edit
$ws.Range("B2").Select()
